$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.593.64"
$ws.Range("E2").Value = "  -1.03%  "
$ws.Range("D3").Value = "1.859.95"
$ws.Range("E3").Value = "  -1.56%  "
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").Value = "'335.07"
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("E6").Value = "  -0.29%  "
$ws.Range("D7").Value = "'0.4647"
$ws.Range("E7").Value = "  -0.92%  "
$ws.Range("D8").Value = "'0.3914"
$ws.Range("E8").Value = "  -0.17%  "
$ws.Range("D9").Value = "'46.00"
$ws.Range("E9").Value = "  -3.76%  "
$ws.Range("D10").Value = "'0.07951"
$ws.Range("E10").Value = "  -1.51%  "
$ws.Range("D11").Value = "'0.9977"
$ws.Range("E11").Value = "  -2.21%  "
$ws.Range("E12").Value = "  -1.61%  "
$ws.Range("D13").Value = "1.874.57"
$ws.Range("E13").Value = "  -0.66%  "
$ws.Range("D14").Value = "'5.928"
$ws.Range("E14").Value = "  -0.60%  "
$ws.Range("D15").Value = "'7.189"
$ws.Range("E15").Value = "  +1.21%  "
$ws.Range("D16").Value = "'1.015"
$ws.Range("E16").Value = "  -0.35%  "
$ws.Range("E17").Value = "  +0.93%  "
$ws.Range("D18").Value = "'0.06730"
$ws.Range("E18").Value = "  -0.71%  "
$ws.Range("E19").Value = "  -0.71%  "
$ws.Range("E20").Value = "  -0.22%  "
$ws.Range("E21").Value = "  -0.31%  "
$ws.Range("D22").Value = "27.599.14"
$ws.Range("E22").Value = "  -1.10%  "
$ws.Range("E23").Value = "  -1.17%  "
$ws.Range("D24").Value = "'10.93"
$ws.Range("E24").Value = "  -0.56%  "
$ws.Range("D25").Value = "'2.304"
$ws.Range("E25").Value = "  -1.59%  "
$ws.Range("D26").Value = "2.089.25"
$ws.Range("E26").Value = "  -1.21%  "
$ws.Range("D27").Value = "'159.54"
$ws.Range("E27").Value = "  -0.18%  "
$ws.Range("D28").Value = "'19.61"
$ws.Range("D29").Value = "'2.141"
$ws.Range("E29").Value = "  +2.73%  "
$ws.Range("D30").Value = "'5.414"
$ws.Range("E30").Value = "  -1.18%  "
$ws.Range("D31").Value = "'121.62"
$ws.Range("E31").Value = "  -0.34%  "
$ws.Range("D32").Value = "'0.9734"
$ws.Range("E32").Value = "  +0.57%  "
$ws.Range("D33").Value = "'0.09432"
$ws.Range("E33").Value = "  -0.58%  "
$ws.Range("D34").Value = "'3.625"
$ws.Range("E34").Value = "  -0.68%  "
$ws.Range("D35").Value = "'5.302"
$ws.Range("E35").Value = "  -0.88%  "
$ws.Range("D36").Value = "'1.337"
$ws.Range("E36").Value = "  -5.14%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "'0.02230"
$ws.Range("E37").Value = "  -1.05%  "
$ws.Range("B38").Value = "Hedera"
$ws.Range("C38").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D38").Value = "'0.06014"
$ws.Range("E38").Value = "  -1.76%  "
$ws.Range("D39").Value = "'8.314"
$ws.Range("E39").Value = "  +3.46%  "
$ws.Range("D40").Value = "'1.193"
$ws.Range("E40").Value = "  -1.95%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "'0.5935"
$ws.Range("E41").Value = "  -0.96%  "
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").Value = "'0.1866"
$ws.Range("E42").Value = "  -1.10%  "
$ws.Range("B43").Value = "Aptos"
$ws.Range("C43").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D43").Value = "'10.29"
$ws.Range("E43").Value = "  +0.09%  "
$ws.Range("B44").Value = "WEMIXTOKEN"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Value = "'1.251"
$ws.Range("E44").Value = "  -1.15%  "
$ws.Range("B45").Value = "Decentraland"
$ws.Range("C45").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D45").Value = "'0.5584"
$ws.Range("E45").Value = "  -1.89%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'12.12"
$ws.Range("E46").Value = "  -0.07%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "'1.918"
$ws.Range("E47").Value = "  -0.56%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "'0.06734"
$ws.Range("E48").Value = "  -2.87%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "'111.34"
$ws.Range("E49").Value = "  -2.28%  "
$ws.Range("B50").Value = "EOS"
$ws.Range("C50").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D50").Value = "'1.050"
$ws.Range("E50").Value = "  -1.97%  "
$ws.Range("B51").Value = "PaxDollar"
$ws.Range("C51").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D51").Value = "'1.015"
$ws.Range("E51").Value = "  -0.23%  "
